# The document ends with a single empty (bold, 22pt) paragraph that carries
# the "_GoBack" bookmark (Word's "last edit location" marker). The edit
# turns that one paragraph into six paragraphs:
#   - the original paragraph, now empty (no bookmark)
#   - three more empty (bold, 22pt) paragraphs
#   - a new (bold, 26pt) heading paragraph reading
#     "YOU BETTER SEE THIS ONE OR ELSE..." immediately followed by the
#     _GoBack bookmark (i.e. the bookmark now marks the spot right after the
#     freshly typed text - exactly where Word leaves it after you type and
#     then press Enter once more)
#   - one final empty (bold, 22pt) paragraph
#
# We locate the paragraph that owns the bookmark (rather than assuming a
# fixed paragraph index), then replace that paragraph's raw OOXML with the
# six-paragraph fragment below via Range.InsertXML, which lets us place the
# bookmark and the new runs/paragraph marks exactly as needed without Word
# leaving stray empty runs behind.

$d = $word.ActiveDocument

$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Range.Start

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $bmStart -and $bmStart -lt $p.Range.End) {
        $target = $p
    }
}

$r = $target.Range

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$emptyPara = "<w:p $wns><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val='44'/><w:szCs w:val='44'/></w:rPr></w:pPr></w:p>"

$headingPara = "<w:p $wns>" +
    "<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val='52'/><w:szCs w:val='52'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val='52'/><w:szCs w:val='52'/></w:rPr>" +
    "<w:t>YOU BETTER SEE THIS ONE OR ELSE&#8230;</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
    "</w:p>"

$xml = ($emptyPara * 4) + $headingPara + $emptyPara

$r.InsertXML($xml)
